# Apply the "Updated symbol list" data refresh to the crypto price sheet.
# All Price (column D) values are stored as text in this workbook, so a
# leading apostrophe is used when assigning numeric-looking strings to
# force Excel to keep them as text (quotePrefix) instead of converting
# them to real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) text-number updates -------------------------------
$ws.Range("D4").Value  = "'6.196"
$ws.Range("D5").Value  = "'0.06104"
$ws.Range("D6").Value  = "'6.739"
$ws.Range("D7").Value  = "'3.500"
$ws.Range("D9").Value  = "'0.7984"
$ws.Range("D10").Value = "'0.1573"
$ws.Range("D11").Value = "'0.08110"
$ws.Range("D12").Value = "'0.03344"
$ws.Range("D13").Value = "'0.03108"
$ws.Range("D14").Value = "'0.09296"
$ws.Range("D15").Value = "'3.922"
$ws.Range("D16").Value = "'0.001687"
$ws.Range("D17").Value = "'0.04820"
$ws.Range("D19").Value = "'0.006235"
$ws.Range("D20").Value = "'0.001097"
$ws.Range("D21").Value = "'0.003398"
$ws.Range("D23").Value = "'3.693"
$ws.Range("D25").Value = "'0.3363"
$ws.Range("D40").Value = "'0.04586"

# --- Rows 41-43 : symbol list was re-ranked (3-way rotation) -----------
# Row 41 becomes what used to be KickToken (row 43), with refreshed price
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.007130"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 becomes what used to be BKEXToken (row 41), with refreshed price
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1122"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 becomes what used to be CEJI (row 42), with refreshed price
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003130"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining column D (Price) text-number updates ---------------------
$ws.Range("D45").Value = "'0.002970"
$ws.Range("D46").Value = "'0.00006000"
$ws.Range("D48").Value = "'0.7499"
$ws.Range("D49").Value = "'0.06105"
$ws.Range("D50").Value = "'0.00002100"
